$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(21, 1).Value = 2
$ws.Cells.Item(21, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(21, 2).Font.Italic = $true
$ws.Cells.Item(21, 3).Value = "TCMNG01"
$ws.Cells.Item(21, 4).Value = "Y"

$ws.Cells.Item(22, 1).Value = 2
$ws.Cells.Item(22, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(22, 2).Font.Italic = $true
$ws.Cells.Item(22, 3).Value = "TCMNG02"
$ws.Cells.Item(22, 4).Value = "Y"

$ws.Cells.Item(23, 1).Value = 2
$ws.Cells.Item(23, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(23, 2).Font.Italic = $true
$ws.Cells.Item(23, 3).Value = "TCMNG03"
$ws.Cells.Item(23, 4).Value = "Y"

$ws.Cells.Item(24, 1).Value = 2
$ws.Cells.Item(24, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(24, 2).Font.Italic = $true
$ws.Cells.Item(24, 3).Value = "TCMNG04"
$ws.Cells.Item(24, 4).Value = "Y"

$ws.Cells.Item(25, 1).Value = 2
$ws.Cells.Item(25, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(25, 2).Font.Italic = $true
$ws.Cells.Item(25, 3).Value = "TCMNG05"
$ws.Cells.Item(25, 4).Value = "Y"

$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(26, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(26, 2).Font.Italic = $true
$ws.Cells.Item(26, 3).Value = "TCMNG06"
$ws.Cells.Item(26, 4).Value = "Y"

$ws.Cells.Item(27, 1).Value = 2
$ws.Cells.Item(27, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(27, 2).Font.Italic = $true
$ws.Cells.Item(27, 3).Value = "TCMNG07"
$ws.Cells.Item(27, 4).Value = "Y"

$ws.Cells.Item(28, 1).Value = 2
$ws.Cells.Item(28, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(28, 2).Font.Italic = $true
$ws.Cells.Item(28, 3).Value = "TCMNG08"
$ws.Cells.Item(28, 4).Value = "Y"

$ws.Cells.Item(29, 1).Value = 2
$ws.Cells.Item(29, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(29, 2).Font.Italic = $true
$ws.Cells.Item(29, 3).Value = "TCMNG09"
$ws.Cells.Item(29, 4).Value = "Y"

$ws.Cells.Item(30, 1).Value = 2
$ws.Cells.Item(30, 2).Value = "Nansenia groenlandica"
$ws.Cells.Item(30, 2).Font.Italic = $true
$ws.Cells.Item(30, 3).Value = "TCMNG10"
$ws.Cells.Item(30, 4).Value = "Y"

$ws.Cells.Item(31, 1).Value = 6
$ws.Cells.Item(31, 2).Value = "Arctozenus risso"
$ws.Cells.Item(31, 2).Font.Italic = $true
$ws.Cells.Item(31, 3).Value = "TCAR001"
$ws.Cells.Item(31, 4).Value = "Y"

$ws.Cells.Item(32, 1).Value = 6
$ws.Cells.Item(32, 2).Value = "Arctozenus risso"
$ws.Cells.Item(32, 2).Font.Italic = $true
$ws.Cells.Item(32, 3).Value = "TCAR002"
$ws.Cells.Item(32, 4).Value = "Y"

$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "Arctozenus risso"
$ws.Cells.Item(33, 2).Font.Italic = $true
$ws.Cells.Item(33, 3).Value = "TCAR003"
$ws.Cells.Item(33, 4).Value = "Y"

$ws.Cells.Item(34, 1).Value = 6
$ws.Cells.Item(34, 2).Value = "Arctozenus risso"
$ws.Cells.Item(34, 2).Font.Italic = $true
$ws.Cells.Item(34, 3).Value = "TCAR004"
$ws.Cells.Item(34, 4).Value = "Y"

$ws.Cells.Item(35, 1).Value = 6
$ws.Cells.Item(35, 2).Value = "Arctozenus risso"
$ws.Cells.Item(35, 2).Font.Italic = $true
$ws.Cells.Item(35, 3).Value = "TCAR005"
$ws.Cells.Item(35, 4).Value = "Y"

$ws.Cells.Item(36, 1).Value = 6
$ws.Cells.Item(36, 2).Value = "Arctozenus risso"
$ws.Cells.Item(36, 2).Font.Italic = $true
$ws.Cells.Item(36, 3).Value = "TCAR006"
$ws.Cells.Item(36, 4).Value = "Y"

$ws.Cells.Item(37, 1).Value = 6
$ws.Cells.Item(37, 2).Value = "Arctozenus risso"
$ws.Cells.Item(37, 2).Font.Italic = $true
$ws.Cells.Item(37, 3).Value = "TCAR007"
$ws.Cells.Item(37, 4).Value = "Y"

$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(38, 2).Value = "Arctozenus risso"
$ws.Cells.Item(38, 2).Font.Italic = $true
$ws.Cells.Item(38, 3).Value = "TCAR008"
$ws.Cells.Item(38, 4).Value = "Y"

$ws.Cells.Item(39, 1).Value = 6
$ws.Cells.Item(39, 2).Value = "Arctozenus risso"
$ws.Cells.Item(39, 2).Font.Italic = $true
$ws.Cells.Item(39, 3).Value = "TCAR009"
$ws.Cells.Item(39, 4).Value = "Y"

$ws.Cells.Item(40, 1).Value = 6
$ws.Cells.Item(40, 2).Value = "Arctozenus risso"
$ws.Cells.Item(40, 2).Font.Italic = $true
$ws.Cells.Item(40, 3).Value = "TCAR010"
$ws.Cells.Item(40, 4).Value = "Y"

$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = "Arctozenus risso"
$ws.Cells.Item(41, 2).Font.Italic = $true
$ws.Cells.Item(41, 3).Value = "TCAR011"
$ws.Cells.Item(41, 4).Value = "Y"

$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = "Arctozenus risso"
$ws.Cells.Item(42, 2).Font.Italic = $true
$ws.Cells.Item(42, 3).Value = "TCAR012"
$ws.Cells.Item(42, 4).Value = "Y"

$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(43, 2).Value = "Arctozenus risso"
$ws.Cells.Item(43, 2).Font.Italic = $true
$ws.Cells.Item(43, 3).Value = "TCAR013"
$ws.Cells.Item(43, 4).Value = "Y"

$ws.Cells.Item(44, 1).Value = 2
$ws.Cells.Item(44, 2).Value = "Arctozenus risso"
$ws.Cells.Item(44, 2).Font.Italic = $true
$ws.Cells.Item(44, 3).Value = "TCAR014"
$ws.Cells.Item(44, 4).Value = "Y"

$ws.Cells.Item(45, 1).Value = 2
$ws.Cells.Item(45, 2).Value = "Arctozenus risso"
$ws.Cells.Item(45, 2).Font.Italic = $true
$ws.Cells.Item(45, 3).Value = "TCAR015"
$ws.Cells.Item(45, 4).Value = "Y"

$ws.Cells.Item(46, 1).Value = 2
$ws.Cells.Item(46, 2).Value = "Arctozenus risso"
$ws.Cells.Item(46, 2).Font.Italic = $true
$ws.Cells.Item(46, 3).Value = "TCAR016"
$ws.Cells.Item(46, 4).Value = "Y"

$ws.Cells.Item(47, 1).Value = 2
$ws.Cells.Item(47, 2).Value = "Arctozenus risso"
$ws.Cells.Item(47, 2).Font.Italic = $true
$ws.Cells.Item(47, 3).Value = "TCAR017"
$ws.Cells.Item(47, 4).Value = "Y"

$ws.Cells.Item(48, 1).Value = 2
$ws.Cells.Item(48, 2).Value = "Arctozenus risso"
$ws.Cells.Item(48, 2).Font.Italic = $true
$ws.Cells.Item(48, 3).Value = "TCAR018"
$ws.Cells.Item(48, 4).Value = "Y"

$ws.Cells.Item(49, 1).Value = 2
$ws.Cells.Item(49, 2).Value = "Arctozenus risso"
$ws.Cells.Item(49, 2).Font.Italic = $true
$ws.Cells.Item(49, 3).Value = "TCAR019"
$ws.Cells.Item(49, 4).Value = "Y"

$ws.Cells.Item(50, 1).Value = 2
$ws.Cells.Item(50, 2).Value = "Arctozenus risso"
$ws.Cells.Item(50, 2).Font.Italic = $true
$ws.Cells.Item(50, 3).Value = "TCAR020"
$ws.Cells.Item(50, 4).Value = "Y"

$ws.Cells.Item(51, 1).Value = 2
$ws.Cells.Item(51, 2).Value = "Arctozenus risso"
$ws.Cells.Item(51, 2).Font.Italic = $true
$ws.Cells.Item(51, 3).Value = "TCAR021"
$ws.Cells.Item(51, 4).Value = "Y"

$ws.Cells.Item(52, 1).Value = 2
$ws.Cells.Item(52, 2).Value = "Arctozenus risso"
$ws.Cells.Item(52, 2).Font.Italic = $true
$ws.Cells.Item(52, 3).Value = "TCAR022"
$ws.Cells.Item(52, 4).Value = "Y"

$ws.Cells.Item(53, 1).Value = 2
$ws.Cells.Item(53, 2).Value = "Arctozenus risso"
$ws.Cells.Item(53, 2).Font.Italic = $true
$ws.Cells.Item(53, 3).Value = "TCAR023"
$ws.Cells.Item(53, 4).Value = "Y"

$ws.Cells.Item(54, 1).Value = 2
$ws.Cells.Item(54, 2).Value = "Arctozenus risso"
$ws.Cells.Item(54, 2).Font.Italic = $true
$ws.Cells.Item(54, 3).Value = "TCAR024"
$ws.Cells.Item(54, 4).Value = "Y"

$ws.Cells.Item(55, 1).Value = 2
$ws.Cells.Item(55, 2).Value = "Arctozenus risso"
$ws.Cells.Item(55, 2).Font.Italic = $true
$ws.Cells.Item(55, 3).Value = "TCAR025"
$ws.Cells.Item(55, 4).Value = "Y"

$ws.Cells.Item(56, 1).Value = 2
$ws.Cells.Item(56, 2).Value = "Arctozenus risso"
$ws.Cells.Item(56, 2).Font.Italic = $true
$ws.Cells.Item(56, 3).Value = "TCAR026"
$ws.Cells.Item(56, 4).Value = "Y"

$ws.Cells.Item(57, 1).Value = 2
$ws.Cells.Item(57, 2).Value = "Arctozenus risso"
$ws.Cells.Item(57, 2).Font.Italic = $true
$ws.Cells.Item(57, 3).Value = "TCAR027"
$ws.Cells.Item(57, 4).Value = "Y"

$ws.Cells.Item(58, 1).Value = 2
$ws.Cells.Item(58, 2).Value = "Arctozenus risso"
$ws.Cells.Item(58, 2).Font.Italic = $true
$ws.Cells.Item(58, 3).Value = "TCAR028"
$ws.Cells.Item(58, 4).Value = "Y"

$ws.Cells.Item(59, 1).Value = 2
$ws.Cells.Item(59, 2).Value = "Arctozenus risso"
$ws.Cells.Item(59, 2).Font.Italic = $true
$ws.Cells.Item(59, 3).Value = "TCAR029"
$ws.Cells.Item(59, 4).Value = "Y"

$ws.Cells.Item(60, 1).Value = 2
$ws.Cells.Item(60, 2).Value = "Arctozenus risso"
$ws.Cells.Item(60, 2).Font.Italic = $true
$ws.Cells.Item(60, 3).Value = "TCAR030"
$ws.Cells.Item(60, 4).Value = "Y"

$ws.Cells.Item(61, 1).Value = 2
$ws.Cells.Item(61, 2).Value = "Arctozenus risso"
$ws.Cells.Item(61, 2).Font.Italic = $true
$ws.Cells.Item(61, 3).Value = "TCAR031"
$ws.Cells.Item(61, 4).Value = "Y"

$ws.Cells.Item(62, 1).Value = 2
$ws.Cells.Item(62, 2).Value = "Arctozenus risso"
$ws.Cells.Item(62, 2).Font.Italic = $true
$ws.Cells.Item(62, 3).Value = "TCAR032"
$ws.Cells.Item(62, 4).Value = "Y"

$ws.Cells.Item(63, 1).Value = 2
$ws.Cells.Item(63, 2).Value = "Arctozenus risso"
$ws.Cells.Item(63, 2).Font.Italic = $true
$ws.Cells.Item(63, 3).Value = "TCAR033"
$ws.Cells.Item(63, 4).Value = "Y"

$ws.Range("A31:A43").Select()
$ws.Application.ActiveWindow.ScrollRow = 29
